# Insert a new data row at row 85 (pushing existing rows 85-174 down to 86-175)
# and populate it with the new "Berenjena" price record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(85).Insert()

$ws.Range("A85").Value = 8
$ws.Range("B85").Value = "Terminal La Palmera de La Serena"
$ws.Range("C85").Value = "Coquimbo"
$ws.Range("D85").Value = 44874
$ws.Range("E85").Value = 4
$ws.Range("F85").Value = 100112001
$ws.Range("G85").Value = "Berenjena"
$ws.Range("H85").Value = "Sin especificar"
$ws.Range("I85").Value = "Primera"
$ws.Range("J85").Value = 530
$ws.Range("K85").Value = 11000
$ws.Range("L85").Value = 12000
$ws.Range("M85").Value = 11500
$ws.Range("N85").Value = "$/caja 40 unidades"
$ws.Range("O85").Value = "Región de Arica y Parinacota"
$ws.Range("P85").Value = 288
$ws.Range("Q85").Value = 40
$ws.Range("R85").Value = "Hortaliza"
